$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 26000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 26000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 26000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -26224

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 26000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 26000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 26000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -26382

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1790.3507
$ws.Range("J17").Value = 1790.3507
$ws.Range("L17").Value = 5371.0521
$ws.Range("N17").Value = -5707.0521

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 33703.207
$ws.Range("I106").Value = 1010
$ws.Range("J106").Value = 197169.25
$ws.Range("K106").Value = 1010
$ws.Range("L106").Value = 197169.25
$ws.Range("M106").Value = -379
$ws.Range("N106").Value = -198431.25

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 20762.184
$ws.Range("I132").Value = 2889.7104
$ws.Range("J132").Value = 82503.45
$ws.Range("K132").Value = 8669.1312
$ws.Range("L132").Value = 247510.35
$ws.Range("M132").Value = -6139.1312
$ws.Range("N132").Value = -252570.35

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3087.3809
$ws.Range("J141").Value = 6369
$ws.Range("L141").Value = 19107
$ws.Range("N141").Value = -29467

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10184.984
$ws.Range("I32").Value = 9583
$ws.Range("K32").Value = 9583
$ws.Range("M32").Value = -9296

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2038.8049
$ws.Range("I74").Value = 1697.75
$ws.Range("J74").Value = 3251.4443
$ws.Range("K74").Value = 1697.75
$ws.Range("L74").Value = 3251.4443
$ws.Range("M74").Value = -823.75
$ws.Range("N74").Value = -4999.4443

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2038.8049
$ws.Range("I77").Value = 1697.75
$ws.Range("J77").Value = 3251.4443
$ws.Range("K77").Value = 8488.75
$ws.Range("L77").Value = 16257.2215
$ws.Range("M77").Value = -4120.75
$ws.Range("N77").Value = -24993.2215

# ARM row 114
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# ARM row 119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 111159880
$ws.Range("J119").Value = 111159880
$ws.Range("L119").Value = 111159880
$ws.Range("N119").Value = -111169556

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 54483.5
$ws.Range("J135").Value = 54483.5
$ws.Range("L135").Value = 54483.5
$ws.Range("N135").Value = -64623.5

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 45689.926
$ws.Range("J139").Value = 45689.926
$ws.Range("L139").Value = 45689.926
$ws.Range("N139").Value = -55969.926

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3968.25
$ws.Range("I107").Value = 2955.5
$ws.Range("J107").Value = 7006.5
$ws.Range("K107").Value = 2955.5
$ws.Range("L107").Value = 7006.5
$ws.Range("M107").Value = -1035.5
$ws.Range("N107").Value = -10846.5

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 34136.445
$ws.Range("J135").Value = 34136.445
$ws.Range("L135").Value = 34136.445
$ws.Range("N135").Value = -44276.445

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3114.09
$ws.Range("I31").Value = 1597.1464
$ws.Range("J31").Value = 4168.2373
$ws.Range("K31").Value = 1597.1464
$ws.Range("L31").Value = 4168.2373
$ws.Range("M31").Value = -1302.1464
$ws.Range("N31").Value = -4758.2373

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3114.09
$ws.Range("I34").Value = 1597.1464
$ws.Range("J34").Value = 4168.2373
$ws.Range("K34").Value = 1597.1464
$ws.Range("L34").Value = 4168.2373
$ws.Range("M34").Value = -1395.1464
$ws.Range("N34").Value = -4572.2373

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 43540.668
$ws.Range("J95").Value = 43540.668
$ws.Range("L95").Value = 43540.668
$ws.Range("N95").Value = -49032.668

# CRP row 119
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 58758.375
$ws.Range("J119").Value = 58758.375
$ws.Range("L119").Value = 58758.375
$ws.Range("N119").Value = -68434.375

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 878917.75
$ws.Range("I132").Value = 2357.2
$ws.Range("J132").Value = 2339852
$ws.Range("K132").Value = 7071.599999999999
$ws.Range("L132").Value = 7019556
$ws.Range("M132").Value = -4541.599999999999
$ws.Range("N132").Value = -7024616

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1363.9398
$ws.Range("I68").Value = 1044.579
$ws.Range("J68").Value = 1458.75
$ws.Range("K68").Value = 3133.737
$ws.Range("L68").Value = 4376.25
$ws.Range("M68").Value = -2322.737
$ws.Range("N68").Value = -5998.25

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1363.9398
$ws.Range("I71").Value = 1044.579
$ws.Range("J71").Value = 1458.75
$ws.Range("K71").Value = 9401.210999999999
$ws.Range("L71").Value = 13128.75
$ws.Range("M71").Value = -5345.210999999999
$ws.Range("N71").Value = -21240.75

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 104925.555
$ws.Range("J129").Value = 1720.238
$ws.Range("L129").Value = 5160.714
$ws.Range("N129").Value = -15160.714

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1004.6667
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 1007
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 1007
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -4251

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1877.9166
$ws.Range("I113").Value = 1915.25
$ws.Range("J113").Value = 1803.25
$ws.Range("K113").Value = 1915.25
$ws.Range("L113").Value = 1803.25
$ws.Range("M113").Value = 254.75
$ws.Range("N113").Value = -6143.25

# GSM row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 59998.5
$ws.Range("J114").Value = 59998.5
$ws.Range("L114").Value = 59998.5
$ws.Range("N114").Value = -68676.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 957.2857
$ws.Range("I93").Value = 703
$ws.Range("J93").Value = 976.8461
$ws.Range("K93").Value = 703
$ws.Range("L93").Value = 976.8461
$ws.Range("M93").Value = 545
$ws.Range("N93").Value = -3472.8461

# LTW row 119
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 1000000000
$ws.Range("J119").Value = 1000000000
$ws.Range("L119").Value = 1000000000
$ws.Range("N119").Value = -1000009676

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 52492.9
$ws.Range("I122").Value = 64909.875
$ws.Range("J122").Value = 2825
$ws.Range("K122").Value = 194729.625
$ws.Range("L122").Value = 8475
$ws.Range("M122").Value = -192279.625
$ws.Range("N122").Value = -13375

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 100000000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# WVR row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 142895710
$ws.Range("J119").Value = 142895710
$ws.Range("L119").Value = 142895710
$ws.Range("N119").Value = -142905386
